$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'66.929.47"
$ws.Range("E2").Formula = "'  +1.91%  "
$ws.Range("D3").Formula = "'3.111.05"
$ws.Range("E3").Formula = "'  +5.39%  "
$ws.Range("E4").Formula = "'  +0.05%  "
$ws.Range("D5").Formula = "'581.06"
$ws.Range("D6").Formula = "'173.25"
$ws.Range("E6").Formula = "'  +6.77%  "
$ws.Range("D7").Formula = "'0.998"
$ws.Range("E7").Formula = "'  -0.14%  "
$ws.Range("D8").Formula = "'3.105.93"
$ws.Range("E8").Formula = "'  +5.32%  "
$ws.Range("E9").Formula = "'  +1.38%  "
$ws.Range("D10").Formula = "'6.50"
$ws.Range("E10").Formula = "'  -3.77%  "
$ws.Range("E11").Formula = "'  +3.37%  "
$ws.Range("D12").Formula = "'0.484"
$ws.Range("E12").Formula = "'  +5.31%  "
$ws.Range("E13").Formula = "'  +1.55%  "
$ws.Range("D14").Formula = "'37.40"
$ws.Range("E14").Formula = "'  +7.76%  "
$ws.Range("D16").Formula = "'3.626.10"
$ws.Range("E16").Formula = "'  +5.37%  "
$ws.Range("D17").Formula = "'66.906.95"
$ws.Range("E17").Formula = "'  +2.05%  "
$ws.Range("D18").Formula = "'7.21"
$ws.Range("E18").Formula = "'  +2.69%  "
$ws.Range("D19").Formula = "'3.110.48"
$ws.Range("E19").Formula = "'  +5.38%  "
$ws.Range("D20").Formula = "'16.15"
$ws.Range("E20").Formula = "'  +1.35%  "
$ws.Range("D21").Formula = "'485.19"
$ws.Range("E21").Formula = "'  +8.59%  "
$ws.Range("E22").Formula = "'  +3.12%  "
$ws.Range("D23").Formula = "'7.54"
$ws.Range("E23").Formula = "'  +3.27%  "
$ws.Range("D24").Formula = "'84.20"
$ws.Range("E24").Formula = "'  +2.37%  "
$ws.Range("E25").Formula = "'  +5.62%  "
$ws.Range("D26").Formula = "'13.21"
$ws.Range("E26").Formula = "'  +7.46%  "
$ws.Range("D27").Formula = "'10.07"
$ws.Range("E27").Formula = "'  +0.28%  "
$ws.Range("E28").Formula = "'  -0.09%  "
$ws.Range("E29").Formula = "'  -1.44%  "
$ws.Range("E30").Formula = "'  -4.81%  "
$ws.Range("E31").Formula = "'  +3.83%  "
$ws.Range("D32").Formula = "'29.02"
$ws.Range("E32").Formula = "'  +6.36%  "
$ws.Range("E33").Formula = "'  -2.09%  "
$ws.Range("E34").Formula = "'  +2.19%  "
$ws.Range("E35").Formula = "'  +0.07%  "
$ws.Range("B36").Formula = "'Filecoin"
$ws.Range("C36").Formula = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Formula = "'5.92"
$ws.Range("E36").Formula = "'  +3.19%  "
$ws.Range("B37").Formula = "'Mantle"
$ws.Range("C37").Formula = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").Formula = "'1.00"
$ws.Range("E37").Formula = "'  +2.99%  "
$ws.Range("D38").Formula = "'47.89"
$ws.Range("E38").Formula = "'  +5.38%  "
$ws.Range("D39").Formula = "'2.13"
$ws.Range("E39").Formula = "'  +6.69%  "
$ws.Range("D40").Formula = "'50.23"
$ws.Range("E40").Formula = "'  +1.95%  "
$ws.Range("E41").Formula = "'  +4.73%  "
$ws.Range("E42").Formula = "'  +0.33%  "
$ws.Range("E43").Formula = "'  +1.35%  "
$ws.Range("E44").Formula = "'  -1.80%  "
$ws.Range("D45").Formula = "'0.0362"
$ws.Range("E45").Formula = "'  +2.77%  "
$ws.Range("D46").Formula = "'2.837.81"
$ws.Range("E46").Formula = "'  +5.74%  "
$ws.Range("D47").Formula = "'382.24"
$ws.Range("E47").Formula = "'  -1.34%  "
$ws.Range("D48").Formula = "'135.43"
$ws.Range("E48").Formula = "'  +1.38%  "
$ws.Range("E49").Formula = "'  +0.02%  "
$ws.Range("E50").Formula = "'  +4.99%  "
$ws.Range("E51").Formula = "'  +2.69%  "
